# priors_specifications.xlsx - apply "prior computation + beginning posterior" edit:
#  1. Fix two spelling typos in the country names ("Pery" -> "Peru",
#     "United States of Amercia" -> "United States of America").
#  2. Sort the data range A1:D34 alphabetically (ascending) by country (col A).
#  3. Turn on AutoFilter for the header row.
#  4. Re-select cell A3 (what was left selected after the sort/filter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix typos -----------------------------------------------------
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Do "Amercia" first so the newly-appended "Peru" string ends up last in
# the shared-strings table (matches how the edit was actually authored).
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq "United States of Amercia") {
        $cell.Value = "United States of America"
    }
}
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq "Pery") {
        $cell.Value = "Peru"
    }
}

# --- 2. Sort A1:D34 by country (column A), ascending, header row included
$dataRange = $ws.Range("A1:D34")
$sortKeyRange = $ws.Range("A2:A34")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKeyRange, 0, 1, 0, 0)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- 3. AutoFilter on the header row -----------------------------------
$headerRange = $ws.Range("A1:D1")
[void]$headerRange.AutoFilter()

$filterFormula = "=" + $ws.Name + "!`$A`$1:`$D`$1"
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $filterFormula)
$filterName.Visible = $false

# --- 4. Restore the selection to A3 ------------------------------------
[void]$ws.Range("A3").Select()
